# Belmans_Jef_Game.xlsx - "Created health class & refactored avatar class"
#
# The underlying XML diff for this commit boils down to two real data/UI
# changes on the "Basic Game rubric" sheet (the rest of the diff is
# incidental save metadata - fileVersion/build numbers, revisionPtr GUIDs,
# and a SharePoint customXml part re-numbering - that Excel regenerates on
# every save and isn't driven by worksheet edits):
#
#   1. The "Camera" rubric score (B2) goes from 2 -> 3.
#   2. The selection in the frozen bottom-right pane moves from B6 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Game rubric")

# Make sure this sheet is the active one (it already is tabSelected in the
# source file) before touching its selection.
$ws.Activate()

# 1. Update the Camera score from 2 to 3.
$ws.Range("B2").Value = 3

# 2. Move the active cell / selection to B3.
$ws.Range("B3").Select()
